# Char Arrays Programs Added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39: blank separator row, matching the format of the other
# week-separator rows (row 12 / row 28 use the same shaded "blank row" style).
$ws.Range("A12:G12").Copy()
$ws.Range("A39:G39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 40: Week 5 -> "Char Arrays & Strings" topic, first question "Length"
$ws.Cells.Item(40, 1).Value = 5
$ws.Cells.Item(40, 2).Value = "Char Arrays & Strings"
$ws.Cells.Item(40, 4).Value = "Length"

# Row 41: second question under the same topic, "Revere"
$ws.Cells.Item(41, 4).Value = "Revere"

# Move the cursor / scroll position to where the author left off editing
$ws.Activate()
$ws.Range("D43").Select()
$excel.ActiveWindow.ScrollRow = 20
